# Edit: move TSRT sheet to the end (after TBL01), and populate two days of
# fresh monitoring data onto the VRKT and TSRT sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Reorder worksheets: move TSRT after TBL01 -----------------------
$tsrt = $wb.Worksheets.Item("TSRT")
$tbl01 = $wb.Worksheets.Item("TBL01")
$tsrt.Move($null, $tbl01)

# --- helper: write a row of values, keeping text-looking numbers as text
function Set-DataRow {
    param(
        $ws,
        [int]$row,
        [object[]]$values
    )

    # Column A is a plain numeric day-of-month value.
    $ws.Cells.Item($row, 1).Value = $values[0]

    # Columns B..G hold text (even when they look like numbers, e.g. "15.90"),
    # stored in the shared string table. Assigning a quoted formula and then
    # freezing it to a value (copy / paste-values) keeps the text intact
    # without Excel auto-converting it to a number, and without touching the
    # cell's existing style.
    for ($i = 1; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        $text = [string]$values[$i]
        $cell.Formula = '="' + $text + '"'
    }

    $rng = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 7))
    $rng.Copy()
    $rng.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# --- 2. VRKT: add rows for days 29 and 30 --------------------------------
$vrkt = $wb.Worksheets.Item("VRKT")
Set-DataRow $vrkt 2 @(29, "15.90", "3.10", "*", "*", "*", "*")
Set-DataRow $vrkt 3 @(30, "62.58", "4.88", "26.92", "33.05", "7.38", "0.41")

# --- 3. TSRT: add rows for days 25-28 ------------------------------------
$tsrtWs = $wb.Worksheets.Item("TSRT")
Set-DataRow $tsrtWs 2 @(25, "241.92", "243.42", "*", "385.36", "*", "32.58")
Set-DataRow $tsrtWs 3 @(26, "50.62", "110.21", "141.43", "129.32", "92.00", "0.09")
Set-DataRow $tsrtWs 4 @(27, "48.17", "104.12", "127.04", "115.67", "93.50", "0.09")
Set-DataRow $tsrtWs 5 @(28, "49.83", "107.13", "121.81", "118.97", "93.62", "0.07")
